$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran advanced text processing reshuffled which word lands on which
# frequency rank. Row 1's header ("Counts") lives in B1 and is untouched;
# column B counts (the numeric frequencies) stay put, only the word in
# column A for each row is updated to match the new ranking.
$words = @(
  "хлеб",
  "вино",
  "скот",
  "холст",
  "кожа",
  "мед",
  "пиво",
  "сукно",
  "овчина",
  "лошадь",
  "воск",
  "масло",
  "сало",
  "железо",
  "полотно",
  "колеса",
  "Крымскую соль",
  "парча",
  "сено",
  "говядина",
  "позумент",
  "табак",
  "шелк",
  "выбойка",
  "сахар",
  "чулок",
  "лес",
  "лыко",
  "сани",
  "коса",
  "сапог",
  "китайка",
  "ладан",
  "платок",
  "гвоздь",
  "ром",
  "овца",
  "обод",
  "рогожа",
  "замок",
  "веревка",
  "конь",
  "горшок",
  "покроми",
  "котел",
  "роза",
  "дуга",
  "брусья",
  "скотский кожа",
  "бечева",
  "гумми",
  "сковорода",
  "сосуд",
  "хомут",
  "нитка"
)

for ($i = 0; $i -lt $words.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $words[$i]
}
